$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 46 (hunk 1)
$ws.Range("H46").Value = 6750
$ws.Range("I46").Value = 6000
$ws.Range("J46").Value = 7500
$ws.Range("K46").Value = 18000
$ws.Range("L46").Value = 22500
$ws.Range("M46").Value = -17881
$ws.Range("N46").Value = -22738

# Row 60 (hunk 2)
$ws.Range("H60").Value = 6750
$ws.Range("I60").Value = 6000
$ws.Range("J60").Value = 7500
$ws.Range("K60").Value = 18000
$ws.Range("L60").Value = 22500
$ws.Range("M60").Value = -17516
$ws.Range("N60").Value = -23468

# Row 62 (hunk 3)
$ws.Range("H62").Value = 6549678.5
$ws.Range("I62").Value = 9264628
$ws.Range("K62").Value = 9264628
$ws.Range("M62").Value = -9264004

# Row 65 (hunk 4)
$ws.Range("H65").Value = 6549678.5
$ws.Range("I65").Value = 9264628
$ws.Range("K65").Value = 46323140
$ws.Range("M65").Value = -46320020

# Row 137 (hunk 5)
$ws.Range("H137").Value = 22729164
$ws.Range("I137").Value = 40001108
$ws.Range("J137").Value = 2921.2104
$ws.Range("K137").Value = 120003324
$ws.Range("L137").Value = 8763.6312
$ws.Range("M137").Value = -120000774
$ws.Range("N137").Value = -13863.6312

$ws = $wb.Worksheets.Item("ARM")
# Row 26 (hunk 6)
$ws.Range("H26").Value = 930.75
$ws.Range("I26").Value = 778
$ws.Range("J26").Value = 2000
$ws.Range("K26").Value = 778
$ws.Range("L26").Value = 2000
$ws.Range("M26").Value = -448
$ws.Range("N26").Value = -2660

# Row 32 (hunk 7)
$ws.Range("H32").Value = 19727.877
$ws.Range("I32").Value = 5258.6196
$ws.Range("J32").Value = 122459.6
$ws.Range("K32").Value = 5258.6196
$ws.Range("L32").Value = 122459.6
$ws.Range("M32").Value = -4971.6196
$ws.Range("N32").Value = -123033.6

# Row 36 (hunk 8)
$ws.Range("H36").Value = 16250
$ws.Range("I36").Value = 16250
$ws.Range("K36").Value = 16250
$ws.Range("M36").Value = -15904

# Row 122 (hunk 9)
$ws.Range("H122").Value = 2137.2
$ws.Range("I122").Value = 1943.5
$ws.Range("J122").Value = 2266.3333
$ws.Range("K122").Value = 5830.5
$ws.Range("L122").Value = 6798.999899999999
$ws.Range("M122").Value = -3380.5
$ws.Range("N122").Value = -11698.9999

# Row 132 (hunk 10)
$ws.Range("H132").Value = 2487.465
$ws.Range("I132").Value = 2014.8379
$ws.Range("J132").Value = 5402
$ws.Range("K132").Value = 6044.5137
$ws.Range("L132").Value = 16206
$ws.Range("M132").Value = -3514.5137
$ws.Range("N132").Value = -21266

# Row 133 (hunk 11)
$ws.Range("H133").Value = 50099.8
$ws.Range("J133").Value = 50099.8
$ws.Range("L133").Value = 50099.8
$ws.Range("N133").Value = -55159.8

$ws = $wb.Worksheets.Item("BSM")
# Row 29 (hunk 12)
$ws.Range("H29").Value = 17244
$ws.Range("I29").Value = 1366
$ws.Range("J29").Value = 49000
$ws.Range("K29").Value = 1366
$ws.Range("L29").Value = 49000
$ws.Range("M29").Value = -1077
$ws.Range("N29").Value = -49578

# Row 118 (hunk 13)
$ws.Range("H118").Value = 26992.5
$ws.Range("J118").Value = 26992.5
$ws.Range("L118").Value = 26992.5
$ws.Range("N118").Value = -30306.5

$ws = $wb.Worksheets.Item("CRP")
# Row 132 (hunk 14)
$ws.Range("H132").Value = 3625024
$ws.Range("I132").Value = 4387252.5
$ws.Range("K132").Value = 13161757.5
$ws.Range("M132").Value = -13159227.5

$ws = $wb.Worksheets.Item("CUL")
# Row 12 (hunk 15)
$ws.Range("H12").Value = 90.888885
$ws.Range("I12").Value = 135.71428
$ws.Range("J12").Value = 62.363636
$ws.Range("K12").Value = 407.14284
$ws.Range("L12").Value = 187.090908
$ws.Range("M12").Value = -234.14284
$ws.Range("N12").Value = -533.090908

# Row 76 (hunk 16)
$ws.Range("H76").Value = 4500
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 13500
$ws.Range("M76").Value = $null
$ws.Range("N76").Value = -14266

# Row 79 (hunk 17)
$ws.Range("H79").Value = 4500
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 13500
$ws.Range("M79").Value = $null
$ws.Range("N79").Value = -16152

# Row 82 (hunk 18)
$ws.Range("H82").Value = 998
$ws.Range("I82").Value = 998
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2994
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2588
$ws.Range("N82").Value = $null

# Row 85 (hunk 19)
$ws.Range("H85").Value = 998
$ws.Range("I85").Value = 998
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2994
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1590
$ws.Range("N85").Value = $null

# Row 100 (hunk 20)
$ws.Range("H100").Value = 5007
$ws.Range("J100").Value = 5007
$ws.Range("L100").Value = 15021
$ws.Range("N100").Value = -16643

# Row 113 (hunk 21)
$ws.Range("H113").Value = 3953318.8
$ws.Range("I113").Value = 661
$ws.Range("J113").Value = 5348374.5
$ws.Range("K113").Value = 1983
$ws.Range("L113").Value = 16045123.5
$ws.Range("M113").Value = 187
$ws.Range("N113").Value = -16049463.5

# Row 131 (hunk 22)
$ws.Range("H131").Value = 6668259.5
$ws.Range("I131").Value = 660
$ws.Range("J131").Value = 7753682.5
$ws.Range("K131").Value = 1980
$ws.Range("L131").Value = 23261047.5
$ws.Range("M131").Value = 3060
$ws.Range("N131").Value = -23271127.5

$ws = $wb.Worksheets.Item("GSM")
# Row 57 (hunk 23)
$ws.Range("H57").Value = 30000
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = $null

# Row 70 (hunk 24)
$ws.Range("H70").Value = 4947.7236
$ws.Range("I70").Value = 5018.2163
$ws.Range("J70").Value = 4686.9
$ws.Range("K70").Value = 5018.2163
$ws.Range("L70").Value = 4686.9
$ws.Range("M70").Value = -4748.2163
$ws.Range("N70").Value = -5226.9

# Row 73 (hunk 25)
$ws.Range("H73").Value = 4947.7236
$ws.Range("I73").Value = 5018.2163
$ws.Range("J73").Value = 4686.9
$ws.Range("K73").Value = 5018.2163
$ws.Range("L73").Value = 4686.9
$ws.Range("M73").Value = -4082.2163
$ws.Range("N73").Value = -6558.9

# Row 80 (hunk 26)
$ws.Range("H80").Value = 33336030
$ws.Range("I80").Value = 2550.35
$ws.Range("J80").Value = 100002990
$ws.Range("K80").Value = 2550.35
$ws.Range("L80").Value = 100002990
$ws.Range("M80").Value = -1552.35
$ws.Range("N80").Value = -100004986

# Row 83 (hunk 27)
$ws.Range("H83").Value = 33336030
$ws.Range("I83").Value = 2550.35
$ws.Range("J83").Value = 100002990
$ws.Range("K83").Value = 12751.75
$ws.Range("L83").Value = 500014950
$ws.Range("M83").Value = -7759.75
$ws.Range("N83").Value = -500024934

# Row 132 (hunk 28)
$ws.Range("H132").Value = 2603.7073
$ws.Range("I132").Value = 2319.8215
$ws.Range("J132").Value = 3215.1538
$ws.Range("K132").Value = 6959.4645
$ws.Range("L132").Value = 9645.4614
$ws.Range("M132").Value = -4429.4645
$ws.Range("N132").Value = -14705.4614

# Row 137 (hunk 29)
$ws.Range("H137").Value = 60000
$ws.Range("J137").Value = 60000
$ws.Range("L137").Value = 60000
$ws.Range("N137").Value = -70200

# Row 138 (hunk 30)
$ws.Range("H138").Value = 79619
$ws.Range("J138").Value = 79619
$ws.Range("L138").Value = 79619
$ws.Range("N138").Value = -89899

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (hunk 31)
$ws.Range("H22").Value = 438.76923
$ws.Range("I22").Value = 463.81818
$ws.Range("J22").Value = 301
$ws.Range("K22").Value = 463.81818
$ws.Range("L22").Value = 301
$ws.Range("M22").Value = -168.81818
$ws.Range("N22").Value = -891

# Row 27 (hunk 32)
$ws.Range("H27").Value = 438.76923
$ws.Range("I27").Value = 463.81818
$ws.Range("J27").Value = 301
$ws.Range("K27").Value = 463.81818
$ws.Range("L27").Value = 301
$ws.Range("M27").Value = -356.81818
$ws.Range("N27").Value = -515

# Row 55 (hunk 33)
$ws.Range("H55").Value = 233.36
$ws.Range("I55").Value = 218.21053
$ws.Range("J55").Value = 281.33334
$ws.Range("K55").Value = 218.21053
$ws.Range("L55").Value = 281.33334
$ws.Range("M55").Value = -45.21053000000001
$ws.Range("N55").Value = -627.33334

# Row 132 (hunk 34)
$ws.Range("H132").Value = 3998.1428
$ws.Range("I132").Value = 2408.7856
$ws.Range("J132").Value = 7176.857
$ws.Range("K132").Value = 7226.3568
$ws.Range("L132").Value = 21530.571
$ws.Range("M132").Value = -4696.3568
$ws.Range("N132").Value = -26590.571

# Row 136 (hunk 35)
$ws.Range("H136").Value = 5426.125
$ws.Range("I136").Value = 2127.1333
$ws.Range("J136").Value = 10924.444
$ws.Range("K136").Value = 6381.3999
$ws.Range("L136").Value = 32773.33199999999
$ws.Range("M136").Value = -3831.3999
$ws.Range("N136").Value = -37873.33199999999

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (hunk 36)
$ws.Range("H132").Value = 3263.1064
$ws.Range("I132").Value = 3122.2424
$ws.Range("J132").Value = 3595.1428
$ws.Range("K132").Value = 9366.727200000001
$ws.Range("L132").Value = 10785.4284
$ws.Range("M132").Value = -6836.727200000001
$ws.Range("N132").Value = -15845.4284

# Row 138 (hunk 37)
$ws.Range("H138").Value = 68823.2
$ws.Range("J138").Value = 68823.2
$ws.Range("L138").Value = 68823.2
$ws.Range("N138").Value = -79103.2
